$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the "Corrected Price" header for the already-present D column data ---
$ws.Range("D1").Value = "Corrected Price"

# --- Update the active selection to D1 (matches the saved selection in the file) ---
$ws.Range("D1").Select()

# --- Insert a clustered column chart sourced from the new Corrected Price column ---
$chartObj = $ws.ChartObjects().Add(100, 20, 300, 150)
$chart = $chartObj.Chart
$chart.ChartType = 51
$chart.SeriesCollection(1).Values = $ws.Range("D2:D4")

$chart.HasTitle = $true
$chart.ChartTitle.Text = "New corrected price chart"

$catAxis = $chart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "New Price"

$valAxis = $chart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "Y-Axis"
